$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bus voltage magnitude results for the 380 kV case (gen setpoint 1.05 -> 1.02 p.u.)
# Columns B:F and I:N are each contiguous per row (column G stays at 1, column H is blank)

$bf2 = New-Object 'object[,]' 1,5
$bf2[0,0] = 1.02
$bf2[0,1] = 1.053688336687904
$bf2[0,2] = 1.053224172736799
$bf2[0,3] = 1.067052245689223
$bf2[0,4] = 1.073819758462312
$ws.Range("B2:F2").Value = $bf2

$in2 = New-Object 'object[,]' 1,6
$in2[0,0] = 1.044018555029285
$in2[0,1] = 1.058704366918162
$in2[0,2] = 1.055970732190669
$in2[0,3] = 1.069761144368334
$in2[0,4] = 1.076510608028805
$in2[0,5] = 1.023093616204831
$ws.Range("I2:N2").Value = $in2

$bf3 = New-Object 'object[,]' 1,5
$bf3[0,0] = 1.02
$bf3[0,1] = 1.055204358522632
$bf3[0,2] = 1.054392980942034
$bf3[0,3] = 1.068549286794199
$bf3[0,4] = 1.075366446613404
$ws.Range("B3:F3").Value = $bf3

$in3 = New-Object 'object[,]' 1,6
$in3[0,0] = 1.044442214793061
$in3[0,1] = 1.059868624403078
$in3[0,2] = 1.056951630958195
$in3[0,3] = 1.071072199742939
$in3[0,4] = 1.077872511106802
$in3[0,5] = 1.023504317320371
$ws.Range("I3:N3").Value = $in3

$bf4 = New-Object 'object[,]' 1,5
$bf4[0,0] = 1.02
$bf4[0,1] = 1.056183620603639
$bf4[0,2] = 1.055147570096798
$bf4[0,3] = 1.069516625305981
$bf4[0,4] = 1.076365958885525
$ws.Range("B4:F4").Value = $bf4

$in4 = New-Object 'object[,]' 1,6
$in4[0,0] = 1.044713962707362
$in4[0,1] = 1.060619873260096
$in4[0,2] = 1.057584057811893
$in4[0,3] = 1.071918683719517
$in4[0,4] = 1.078751952879821
$in4[0,5] = 1.023768786451821
$ws.Range("I4:N4").Value = $in4

$bf5 = New-Object 'object[,]' 1,5
$bf5[0,0] = 1.02
$bf5[0,1] = 1.056594902578915
$bf5[0,2] = 1.055464396108367
$bf5[0,3] = 1.069922979573706
$bf5[0,4] = 1.076785850915315
$ws.Range("B5:F5").Value = $bf5

$in5 = New-Object 'object[,]' 1,6
$in5[0,0] = 1.044827636442696
$in5[0,1] = 1.060935201098663
$in5[0,2] = 1.057849389447092
$in5[0,3] = 1.072274108152175
$in5[0,4] = 1.079121246357433
$in5[0,5] = 1.023879664631491
$ws.Range("I5:N5").Value = $in5

$bf6 = New-Object 'object[,]' 1,5
$bf6[0,0] = 1.02
$bf6[0,1] = 1.056663935433021
$bf6[0,2] = 1.05551756910411
$bf6[0,3] = 1.069991190002979
$bf6[0,4] = 1.076856335099658
$ws.Range("B6:F6").Value = $bf6

$in6 = New-Object 'object[,]' 1,6
$in6[0,0] = 1.044846689468672
$in6[0,1] = 1.060988117061183
$in6[0,2] = 1.057893908246988
$in6[0,3] = 1.072333760030445
$in6[0,4] = 1.079183227711841
$in6[0,5] = 1.023898263766725
$ws.Range("I6:N6").Value = $in6

$bf7 = New-Object 'object[,]' 1,5
$bf7[0,0] = 1.02
$bf7[0,1] = 1.056189117733994
$bf7[0,2] = 1.055151805119128
$bf7[0,3] = 1.069522056259218
$bf7[0,4] = 1.076371570683907
$ws.Range("B7:F7").Value = $bf7

$in7 = New-Object 'object[,]' 1,6
$in7[0,0] = 1.044715483854224
$in7[0,1] = 1.060624088631762
$in7[0,2] = 1.057587605303929
$in7[0,3] = 1.071923434628312
$in7[0,4] = 1.078756889053842
$in7[0,5] = 1.023770269205725
$ws.Range("I7:N7").Value = $in7

$bf8 = New-Object 'object[,]' 1,5
$bf8[0,0] = 1.02
$bf8[0,1] = 1.054201041600414
$bf8[0,2] = 1.053619533963751
$bf8[0,3] = 1.067558461318263
$bf8[0,4] = 1.074342742577213
$ws.Range("B8:F8").Value = $bf8

$in8 = New-Object 'object[,]' 1,6
$in8[0,0] = 1.044162228954423
$in8[0,1] = 1.059098272625486
$in8[0,2] = 1.056302706849705
$in8[0,3] = 1.070204610315209
$in8[0,4] = 1.076971246733027
$in8[0,5] = 1.023232681370149
$ws.Range("I8:N8").Value = $in8

$bf9 = New-Object 'object[,]' 1,5
$bf9[0,0] = 1.02
$bf9[0,1] = 1.050684354862893
$bf9[0,2] = 1.050906143184571
$bf9[0,3] = 1.064087693037093
$bf9[0,4] = 1.070757388537441
$ws.Range("B9:F9").Value = $bf9

$in9 = New-Object 'object[,]' 1,6
$in9[0,0] = 1.043168911161436
$in9[0,1] = 1.056393189943942
$in9[0,2] = 1.054020849829728
$in9[0,3] = 1.067161287600447
$in9[0,4] = 1.073810601644472
$in9[0,5] = 1.022275463068776
$ws.Range("I9:N9").Value = $in9

$bf10 = New-Object 'object[,]' 1,5
$bf10[0,0] = 1.02
$bf10[0,1] = 1.048330323989522
$bf10[0,2] = 1.049087895379035
$bf10[0,3] = 1.061766193076434
$bf10[0,4] = 1.068359727875051
$ws.Range("B10:F10").Value = $bf10

$in10 = New-Object 'object[,]' 1,6
$in10[0,0] = 1.042494149757601
$in10[0,1] = 1.054578368808441
$in10[0,2] = 1.052487375052007
$in10[0,3] = 1.065122189376543
$in10[0,4] = 1.071693558035738
$in10[0,5] = 1.021630504976667
$ws.Range("I10:N10").Value = $in10

$bf11 = New-Object 'object[,]' 1,5
$bf11[0,0] = 1.02
$bf11[0,1] = 1.047308612363881
$bf11[0,2] = 1.0482982832074
$bf11[0,3] = 1.060759033776252
$bf11[0,4] = 1.067319641404722
$ws.Range("B11:F11").Value = $bf11

$in11 = New-Object 'object[,]' 1,6
$in11[0,0] = 1.042198955584949
$in11[0,1] = 1.053789728467163
$in11[0,2] = 1.051820385837317
$in11[0,3] = 1.064236714620773
$in11[0,4] = 1.070774393186369
$in11[0,5] = 1.021349583046477
$ws.Range("I11:N11").Value = $in11

$bf12 = New-Object 'object[,]' 1,5
$bf12[0,0] = 1.02
$bf12[0,1] = 1.046928732103372
$bf12[0,2] = 1.0480046337129
$bf12[0,3] = 1.060384629887008
$bf12[0,4] = 1.066933014266863
$ws.Range("B12:F12").Value = $bf12

$in12 = New-Object 'object[,]' 1,6
$in12[0,0] = 1.042088850510707
$in12[0,1] = 1.053496362102849
$in12[0,2] = 1.051572181371409
$in12[0,3] = 1.063907420839542
$in12[0,4] = 1.070432594290816
$in12[0,5] = 1.021244985306196
$ws.Range("I12:N12").Value = $in12

$bf13 = New-Object 'object[,]' 1,5
$bf13[0,0] = 1.02
$bf13[0,1] = 1.047010234668252
$bf13[0,2] = 1.048067638585096
$bf13[0,3] = 1.060464954524006
$bf13[0,4] = 1.067015960503285
$ws.Range("B13:F13").Value = $bf13

$in13 = New-Object 'object[,]' 1,6
$in13[0,0] = 1.042112489149511
$in13[0,1] = 1.053559309780656
$in13[0,2] = 1.051625442766774
$in13[0,3] = 1.063978073217761
$in13[0,4] = 1.070505928660317
$in13[0,5] = 1.02126743327238
$ws.Range("I13:N13").Value = $in13

$bf14 = New-Object 'object[,]' 1,5
$bf14[0,0] = 1.02
$bf14[0,1] = 1.047277218973182
$bf14[0,2] = 1.04827401727766
$bf14[0,3] = 1.060728091610459
$bf14[0,4] = 1.067287688712707
$ws.Range("B14:F14").Value = $bf14

$in14 = New-Object 'object[,]' 1,6
$in14[0,0] = 1.04218986360378
$in14[0,1] = 1.053765487534619
$in14[0,2] = 1.051799878498368
$in14[0,3] = 1.064209503080576
$in14[0,4] = 1.070746147789659
$in14[0,5] = 1.021340942097774
$ws.Range("I14:N14").Value = $in14

$bf15 = New-Object 'object[,]' 1,5
$bf15[0,0] = 1.02
$bf15[0,1] = 1.047441667332267
$bf15[0,2] = 1.048401127097778
$bf15[0,3] = 1.060890179006438
$bf15[0,4] = 1.067455070351026
$ws.Range("B15:F15").Value = $bf15

$in15 = New-Object 'object[,]' 1,6
$in15[0,0] = 1.042237475924338
$in15[0,1] = 1.053892463205742
$in15[0,2] = 1.051907293623259
$in15[0,3] = 1.064352042846431
$in15[0,4] = 1.070894104082911
$in15[0,5] = 1.021386199989685
$ws.Range("I15:N15").Value = $in15

$bf16 = New-Object 'object[,]' 1,5
$bf16[0,0] = 1.02
$bf16[0,1] = 1.048398080159848
$bf16[0,2] = 1.049140250307993
$bf16[0,3] = 1.061832993360174
$bf16[0,4] = 1.068428714475442
$ws.Range("B16:F16").Value = $bf16

$in16 = New-Object 'object[,]' 1,6
$in16[0,0] = 1.042513676947592
$in16[0,1] = 1.054630648374251
$in16[0,2] = 1.052531577488566
$in16[0,3] = 1.065180901340746
$in16[0,4] = 1.071754507160144
$in16[0,5] = 1.021649113819732
$ws.Range("I16:N16").Value = $in16

$bf17 = New-Object 'object[,]' 1,5
$bf17[0,0] = 1.02
$bf17[0,1] = 1.048997362608279
$bf17[0,2] = 1.049603262377214
$bf17[0,3] = 1.062423871332043
$bf17[0,4] = 1.069038944582029
$ws.Range("B17:F17").Value = $bf17

$in17 = New-Object 'object[,]' 1,6
$in17[0,0] = 1.042686120156578
$in17[0,1] = 1.055092934312603
$in17[0,2] = 1.052922370594503
$in17[0,3] = 1.065700138457478
$in17[0,4] = 1.072293547459451
$in17[0,5] = 1.02181358877065
$ws.Range("I17:N17").Value = $in17

$bf18 = New-Object 'object[,]' 1,5
$bf18[0,0] = 1.02
$bf18[0,1] = 1.049346682992955
$bf18[0,2] = 1.0498731081523
$bf18[0,3] = 1.062768334447768
$bf18[0,4] = 1.069394700402342
$ws.Range("B18:F18").Value = $bf18

$in18 = New-Object 'object[,]' 1,6
$in18[0,0] = 1.042786412397586
$in18[0,1] = 1.055362307216305
$in18[0,2] = 1.053150026087282
$in18[0,3] = 1.066002756993207
$in18[0,4] = 1.072607722736812
$in18[0,5] = 1.021909365137815
$ws.Range("I18:N18").Value = $in18

$bf19 = New-Object 'object[,]' 1,5
$bf19[0,0] = 1.02
$bf19[0,1] = 1.049465753360239
$bf19[0,2] = 1.049965081256954
$bf19[0,3] = 1.062885756202554
$bf19[0,4] = 1.069515973522546
$ws.Range("B19:F19").Value = $bf19

$in19 = New-Object 'object[,]' 1,6
$in19[0,0] = 1.042820560209127
$in19[0,1] = 1.055454110744933
$in19[0,2] = 1.053227602157826
$in19[0,3] = 1.066105901020236
$in19[0,4] = 1.072714808318812
$in19[0,5] = 1.021941995504259
$ws.Range("I19:N19").Value = $in19

$bf20 = New-Object 'object[,]' 1,5
$bf20[0,0] = 1.02
$bf20[0,1] = 1.048933089234351
$bf20[0,2] = 1.049553608508117
$bf20[0,3] = 1.062360495025802
$bf20[0,4] = 1.068973491466792
$ws.Range("B20:F20").Value = $bf20

$in20 = New-Object 'object[,]' 1,6
$in20[0,0] = 1.042667648752699
$in20[0,1] = 1.055043363459843
$in20[0,2] = 1.052880471979297
$in20[0,3] = 1.065644454497558
$in20[0,4] = 1.072235738198997
$in20[0,5] = 1.021795958640685
$ws.Range("I20:N20").Value = $in20

$bf21 = New-Object 'object[,]' 1,5
$bf21[0,0] = 1.02
$bf21[0,1] = 1.047198609077027
$bf21[0,2] = 1.048213253680382
$bf21[0,3] = 1.06065061265698
$bf21[0,4] = 1.067207679686506
$ws.Range("B21:F21").Value = $bf21

$in21 = New-Object 'object[,]' 1,6
$in21[0,0] = 1.042167091393223
$in21[0,1] = 1.053704785265358
$in21[0,2] = 1.05174852413579
$in21[0,3] = 1.064141363563303
$in21[0,4] = 1.070675419809122
$in21[0,5] = 1.021319302527599
$ws.Range("I21:N21").Value = $in21

$bf22 = New-Object 'object[,]' 1,5
$bf22[0,0] = 1.02
$bf22[0,1] = 1.046105920450395
$bf22[0,2] = 1.047368477147256
$bf22[0,3] = 1.059573800832433
$bf22[0,4] = 1.066095745311604
$ws.Range("B22:F22").Value = $bf22

$in22 = New-Object 'object[,]' 1,6
$in22[0,0] = 1.041849726420446
$in22[0,1] = 1.052860674884668
$in22[0,2] = 1.051034187956304
$in22[0,3] = 1.063194055477602
$in22[0,4] = 1.069692181842595
$in22[0,5] = 1.021018157292891
$ws.Range("I22:N22").Value = $in22

$bf23 = New-Object 'object[,]' 1,5
$bf23[0,0] = 1.02
$bf23[0,1] = 1.046685382850977
$bf23[0,2] = 1.047816505105496
$bf23[0,3] = 1.060144807346355
$bf23[0,4] = 1.066685367020937
$ws.Range("B23:F23").Value = $bf23

$in23 = New-Object 'object[,]' 1,6
$in23[0,0] = 1.042018219423697
$in23[0,1] = 1.05330839275552
$in23[0,2] = 1.051413123122504
$in23[0,3] = 1.063696458090546
$in23[0,4] = 1.07021362678536
$in23[0,5] = 1.021177938714871
$ws.Range("I23:N23").Value = $in23

$bf24 = New-Object 'object[,]' 1,5
$bf24[0,0] = 1.02
$bf24[0,1] = 1.048962132319121
$bf24[0,2] = 1.049576045638348
$bf24[0,3] = 1.062389132625544
$bf24[0,4] = 1.06900306747189
$ws.Range("B24:F24").Value = $bf24

$in24 = New-Object 'object[,]' 1,6
$in24[0,0] = 1.042675996084173
$in24[0,1] = 1.055065763230939
$in24[0,2] = 1.052899405048566
$in24[0,3] = 1.065669616435317
$in24[0,4] = 1.072261860450423
$in24[0,5] = 1.021803925429189
$ws.Range("I24:N24").Value = $in24

$bf25 = New-Object 'object[,]' 1,5
$bf25[0,0] = 1.02
$bf25[0,1] = 1.051595150525509
$bf25[0,2] = 1.051609236261452
$bf25[0,3] = 1.064986283180051
$bf25[0,4] = 1.071685560267760
$ws.Range("B25:F25").Value = $bf25

$in25 = New-Object 'object[,]' 1,6
$in25[0,0] = 1.043427906970109
$in25[0,1] = 1.057094503760869
$in25[0,2] = 1.054612896631418
$in25[0,3] = 1.067949827408631
$in25[0,4] = 1.074629423666934
$in25[0,5] = 1.022524117154669
$ws.Range("I25:N25").Value = $in25

